$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert a new table row before row 505 (shifts 505:574 -> 506:575) ---
$ws.Rows.Item(505).Insert(-4121) | Out-Null   # xlShiftDown

# Copy formatting from the row above (504) into the freshly inserted blank row (505)
$ws.Range("A504:K504").Copy() | Out-Null
$ws.Range("A505:K505").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the "mirror EARNED" formula in G505 (PasteSpecial only copies formats)
$ws.Range("G505").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Grow Table1 so it covers the new row (A8:K574 -> A8:K575)
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K575"))

# The resize above can mangle the calculated-column formula text/cache on the
# brand-new last row (575) - re-assert it explicitly so it stays healthy.
$ws.Range("G575").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# A505 must stay truly blank (no date) - Insert/PasteSpecial above leaves the
# value from the old row504 copy behind as format-only, but make sure it's clear.
$ws.Range("A505").ClearContents()

# --- 2. Fill in the missing monthly VL accruals for May & June 2023 ---
$ws.Range("C502").Value = 1.25
$ws.Range("C503").Value = 1.25

# --- 3. Row 504 (7/1/2023): VL used 3 days on 7/17,20,21/2023 ---
$ws.Range("B504").Value = "VL(3-0-0)"
$ws.Range("D504").Value = 3
$ws.Range("K504").Value = "7/17,20,21/2023"

# --- 4. New row 505: SL used 2 days on 7/12-13/2023 ---
$ws.Range("B505").Value = "SL(2-0-0)"
$ws.Range("H505").Value = 2
$ws.Range("K505").Value = "7/12-13/2023"

# --- 5. Selection bookkeeping to match the authored file ---
$ws.Range("I505").Select()
